$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1769
$ws1.Range("F6").Value = 658
$ws1.Range("F8").Value = 473
$ws1.Range("F9").Value = 4419
$ws1.Range("F14").Value = 1280
$ws1.Range("F16").Value = 1890
$ws1.Range("F17").Value = 2961
$ws1.Range("F18").Value = 1794
$ws1.Range("F21").Value = 167
$ws1.Range("F24").Value = 922
$ws1.Range("F25").Value = 297
$ws1.Range("F26").Value = 27
$ws1.Range("F27").Value = 2295
$ws1.Range("F29").Value = 2368
$ws1.Range("F31").Value = 700
$ws1.Range("F32").Value = 548
$ws1.Range("F34").Value = 880
$ws1.Range("F35").Value = 409
$ws1.Range("F37").Value = 891
$ws1.Range("F38").Value = 1162
$ws1.Range("F40").Value = 669
$ws1.Range("F42").Value = 355
$ws1.Range("F43").Value = 276
$ws1.Range("F44").Value = 3466

# Sheet: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 22

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 658
$ws4.Range("F9").Value = 473
$ws4.Range("F10").Value = 4419
$ws4.Range("F16").Value = 1280
$ws4.Range("F17").Value = 2961
$ws4.Range("F19").Value = 1794
$ws4.Range("F22").Value = 167
$ws4.Range("F27").Value = 922
$ws4.Range("F28").Value = 297
$ws4.Range("F29").Value = 2295
$ws4.Range("F33").Value = 2368
$ws4.Range("F34").Value = 700
$ws4.Range("F35").Value = 548
$ws4.Range("F36").Value = 880
$ws4.Range("F38").Value = 891
$ws4.Range("F39").Value = 1162
$ws4.Range("F40").Value = 669
$ws4.Range("F43").Value = 355
$ws4.Range("F47").Value = 276
$ws4.Range("F48").Value = 3466
